$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 6640.5
$ws.Range("I19").Value = 7965.6665
$ws.Range("K19").Value = 7965.6665
$ws.Range("M19").Value = -7790.6665
# row 86
$ws.Range("H86").Value = 2938
$ws.Range("I86").Value = 2968.7144
$ws.Range("J86").Value = 2866.3333
$ws.Range("K86").Value = 2968.7144
$ws.Range("L86").Value = 2866.3333
$ws.Range("M86").Value = -1845.7144
$ws.Range("N86").Value = -5112.3333
# row 89
$ws.Range("H89").Value = 2938
$ws.Range("I89").Value = 2968.7144
$ws.Range("J89").Value = 2866.3333
$ws.Range("K89").Value = 14843.572
$ws.Range("L89").Value = 14331.6665
$ws.Range("M89").Value = -9227.572
$ws.Range("N89").Value = -25563.6665
# row 125
$ws.Range("H125").Value = 16519.62
$ws.Range("I125").Value = 46954.43
$ws.Range("J125").Value = 1302.2142
$ws.Range("K125").Value = 422589.87
$ws.Range("L125").Value = 11719.9278
$ws.Range("M125").Value = -420129.87
$ws.Range("N125").Value = -16639.9278
# row 137
$ws.Range("H137").Value = 3120.36
$ws.Range("I137").Value = 1138.2325
$ws.Range("K137").Value = 3414.6975
$ws.Range("M137").Value = -864.6975000000002
# row 138
$ws.Range("H138").Value = 1980.1333
$ws.Range("I138").Value = 1307.0625
$ws.Range("K138").Value = 3921.1875
$ws.Range("M138").Value = 1218.8125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 724211.6
$ws.Range("I45").Value = 1445094.9
$ws.Range("K45").Value = 1445094.9
$ws.Range("M45").Value = -1444717.9
# row 110
$ws.Range("H110").Value = 2420.9688
$ws.Range("I110").Value = 2687.6155
$ws.Range("K110").Value = 2687.6155
$ws.Range("M110").Value = -642.6154999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 105
$ws.Range("H105").Value = 5459.1
$ws.Range("I105").Value = 6090.6665
$ws.Range("K105").Value = 6090.6665
$ws.Range("M105").Value = -4343.6665
# row 107
$ws.Range("H107").Value = 21716.52
$ws.Range("I107").Value = 39015.668
$ws.Range("J107").Value = 1408.826
$ws.Range("K107").Value = 39015.668
$ws.Range("L107").Value = 1408.826
$ws.Range("M107").Value = -37095.668
$ws.Range("N107").Value = -5248.826

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 2877.9
$ws.Range("I16").Value = 3097.5
$ws.Range("K16").Value = 3097.5
$ws.Range("M16").Value = -2810.5
# row 69
$ws.Range("H69").Value = 8633.333000000001
$ws.Range("I69").Value = 8633.333000000001
$ws.Range("K69").Value = 8633.333000000001
$ws.Range("M69").Value = -7884.333000000001
# row 72
$ws.Range("H72").Value = 8633.333000000001
$ws.Range("I72").Value = 8633.333000000001
$ws.Range("K72").Value = 25899.999
$ws.Range("M72").Value = -22155.999
# row 105
$ws.Range("H105").Value = 1697.6
$ws.Range("I105").Value = 1384.5
$ws.Range("K105").Value = 1384.5
$ws.Range("M105").Value = 362.5
# row 107
$ws.Range("H107").Value = 2608.3333
$ws.Range("I107").Value = 2594.6667
$ws.Range("K107").Value = 2594.6667
$ws.Range("M107").Value = -674.6667000000002
# row 109
$ws.Range("H109").Value = 73999
$ws.Range("J109").Value = 73999
$ws.Range("L109").Value = 73999
$ws.Range("N109").Value = -76079
# row 113
$ws.Range("H113").Value = 2877.9
$ws.Range("I113").Value = 3097.5
$ws.Range("K113").Value = 3097.5
$ws.Range("M113").Value = -927.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Range("H68").Value = 3582.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3582.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 10748.25
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -12370.25
# row 71
$ws.Range("H71").Value = 3582.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3582.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 32244.75
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -40356.75
# row 128
$ws.Range("H128").Value = 310098.8
$ws.Range("I128").Value = 310098.8
$ws.Range("K128").Value = 930296.3999999999
$ws.Range("M128").Value = -925316.3999999999
# row 132
$ws.Range("H132").Value = 1685.9762
$ws.Range("I132").Value = 777.17645
$ws.Range("J132").Value = 2303.96
$ws.Range("K132").Value = 6994.58805
$ws.Range("L132").Value = 20735.64
$ws.Range("M132").Value = -4464.58805
$ws.Range("N132").Value = -25795.64

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 36
$ws.Range("H36").Value = 8399.25
$ws.Range("I36").Value = 8399.25
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8399.25
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -7914.25
$ws.Range("N36").ClearContents()
# row 80
$ws.Range("H80").Value = 3270.8948
$ws.Range("J80").Value = 4376.8887
$ws.Range("L80").Value = 4376.8887
$ws.Range("N80").Value = -6372.8887
# row 83
$ws.Range("H83").Value = 3270.8948
$ws.Range("J83").Value = 4376.8887
$ws.Range("L83").Value = 21884.4435
$ws.Range("N83").Value = -31868.4435
# row 132
$ws.Range("H132").Value = 2639.6667
$ws.Range("I132").Value = 2439.8333
$ws.Range("K132").Value = 7319.499899999999
$ws.Range("M132").Value = -4789.499899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 4820.65
$ws.Range("I7").Value = 5221.6665
$ws.Range("K7").Value = 5221.6665
$ws.Range("M7").Value = -5109.6665
# row 22
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5590
# row 27
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5214
# row 122
$ws.Range("H122").Value = 2850
$ws.Range("I122").Value = 2357.5833
$ws.Range("K122").Value = 7072.749899999999
$ws.Range("M122").Value = -4622.749899999999
# row 126
$ws.Range("H126").Value = 4820.65
$ws.Range("I126").Value = 5221.6665
$ws.Range("K126").Value = 15664.9995
$ws.Range("M126").Value = -13194.9995
# row 132
$ws.Range("H132").Value = 2457.6511
$ws.Range("I132").Value = 1947.0513
$ws.Range("K132").Value = 5841.1539
$ws.Range("M132").Value = -3311.1539
# row 136
$ws.Range("H136").Value = 2347.5676
$ws.Range("I136").Value = 1802.0968
$ws.Range("K136").Value = 5406.2904
$ws.Range("M136").Value = -2856.2904

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 1325.5714
$ws.Range("I113").Value = 1240
$ws.Range("J113").Value = 1359.8
$ws.Range("K113").Value = 3720
$ws.Range("L113").Value = 4079.4
$ws.Range("M113").Value = -1550
$ws.Range("N113").Value = -8419.4
# row 122
$ws.Range("H122").Value = 2926
$ws.Range("I122").Value = 2416.5715
$ws.Range("J122").Value = 4709
$ws.Range("K122").Value = 7249.7145
$ws.Range("L122").Value = 14127
$ws.Range("M122").Value = -4799.7145
$ws.Range("N122").Value = -19027
# row 126
$ws.Range("H126").Value = 2565.9412
$ws.Range("I126").Value = 1967.7273
$ws.Range("J126").Value = 3662.6667
$ws.Range("K126").Value = 5903.1819
$ws.Range("L126").Value = 10988.0001
$ws.Range("M126").Value = -3433.1819
$ws.Range("N126").Value = -15928.0001
# row 136
$ws.Range("H136").Value = 2478.724
$ws.Range("I136").Value = 2478.724
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7436.172
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4886.172
$ws.Range("N136").ClearContents()
